$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1, matching the style of the existing header cells (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Add time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 10:50:54.132523"
$ws.Range("F3").Value = "2021-10-05 10:50:54.132535"
